$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.709.23"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "2.388.16"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "504.20"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.41"
$ws.Range("E6").Value = "  +2.37%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").Value = "2.391.82"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0973"
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.67"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").Value = "2.810.71"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "56.616.50"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.65"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").Value = "2.381.55"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.18"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "309.00"
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.26"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  -4.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.60"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.384"
$ws.Range("E27").Value = "  +3.74%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.42"
$ws.Range("E29").Value = "  +2.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "176.43"
$ws.Range("E30").Value = "  +1.21%  "
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("E33").Value = "  +1.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.85"
$ws.Range("E34").Value = "  -4.50%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.80"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.80"
$ws.Range("E39").Value = "  +1.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.84"
$ws.Range("E40").Value = "  +2.94%  "
$ws.Range("E41").Value = "  +5.72%  "
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "131.33"
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.84"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.566"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0908"
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "249.55"
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0210"
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("E51").Value = "  +7.51%  "
